$d = $word.ActiveDocument

# --- Change 1: bold the word "player" in the "(10 points): As a player, the price ..." paragraph ---
$p7 = $d.Paragraphs.Item(7).Range
$rng = $d.Range($p7.Start, $p7.End)
$found = $rng.Find.Execute("player", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Bold = 1
}

# --- Change 2: relocate the "_GoBack" bookmark from the start of the document
#     to wrap the word "much " inside the "(10 points): ... how much lemonade ..." paragraph ---
$p8 = $d.Paragraphs.Item(8).Range
$target = $d.Range($p8.Start, $p8.End)
$foundMuch = $target.Find.Execute("much ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMuch) {
    $d.Bookmarks.Add("_GoBack", $target)
}
